# Apply "Trade #93 closed at 2026-02-17 09:09:07" update to live_trading_results workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.61   # Current Capital
$summary.Range("B4").Value = -0.38     # Total P&L $
$summary.Range("B6").Value = 93        # Total Trades
$summary.Range("B8").Value = 39        # Losing Trades
$summary.Range("B9").Value = 40.86     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.61      # Capital
$status.Range("D4").Value = 93         # Trades
$status.Range("E4").Value = -0.38      # P&L $
$status.Range("F4").Value = -0.39      # P&L %
$status.Range("G4").Value = 40.86      # Win Rate %

# ---------------------------------------------------------------------------
# Helper to append the new trade row (#93) to a trade log sheet.
# ---------------------------------------------------------------------------
function Add-Trade93Row($ws) {
    $row = 94

    $ws.Cells.Item($row, 1).Value = 93
    # Force the date-looking text to stay as plain text instead of being
    # auto-converted into an Excel date serial number.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"
    $ws.Cells.Item($row, 3).Value = "09:09:01"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 0.13
    $ws.Cells.Item($row, 7).Value = 0.1
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = -23.0769
    $ws.Cells.Item($row, 10).Value = -0.03
    $ws.Cells.Item($row, 11).Value = 99.61
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.14

    # Restore default (General/Normal) styling now that the text values are
    # safely stored, so no extraneous number formatting lingers on the row.
    $ws.Range("A94:Q94").Style = "Normal"
}

# ---------------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-Trade93Row $allTrades

# ---------------------------------------------------------------------------
# MarketMaking sheet (mirrors All Trades for this strategy)
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade93Row $marketMaking
